# Correcting error in xlsx file: fix typo "floa16" -> "float16" in cell F3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = "float16"

# Update the active selection to F3, matching the saved cursor position.
$ws.Range("F3").Select()
